# Updated cryptos list on Fri Nov  3 06:15:48 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while preserving it as plain text (matches the source
# workbook, where every Price/Volume cell is stored as a text string even when
# it is numeric-looking). A leading apostrophe forces text entry; ClearFormats
# then strips the "Text" quote-prefix number format Excel applies automatically,
# leaving the cell on the default (unstyled) format like the rest of the sheet.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "34.729.37"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.808.80"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "232.07"
$ws.Range("E5").Value = "  +1.09%  "
Set-TextValue "D6" "0.603"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue "D8" "39.28"
$ws.Range("E8").Value = "  -7.27%  "
Set-TextValue "D9" "0.320"
$ws.Range("E9").Value = "  +4.80%  "
Set-TextValue "D10" "0.0682"
$ws.Range("E10").Value = "  -1.50%  "
Set-TextValue "D11" "0.0994"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "2.069.16"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "1.805.73"
$ws.Range("E13").Value = "  -2.08%  "
Set-TextValue "D14" "0.665"
$ws.Range("E14").Value = "  -0.40%  "
Set-TextValue "D15" "10.96"
$ws.Range("E15").Value = "  -3.32%  "
Set-TextValue "D16" "4.58"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "34.681.70"
$ws.Range("E17").Value = "  -1.91%  "
Set-TextValue "D18" "69.69"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("E19").Value = "  -1.47%  "
Set-TextValue "D20" "239.99"
$ws.Range("E20").Value = "  -2.45%  "
Set-TextValue "D21" "11.90"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +2.49%  "
Set-TextValue "D25" "172.77"
$ws.Range("E25").Value = "  +2.18%  "
Set-TextValue "D26" "7.75"
$ws.Range("E26").Value = "  -1.67%  "
Set-TextValue "D27" "17.23"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("E28").Value = "  -1.98%  "
Set-TextValue "D29" "1.54"
$ws.Range("E29").Value = "  +11.92%  "
$ws.Range("E30").Value = "  +0.06%  "
Set-TextValue "D31" "4.04"
$ws.Range("E31").Value = "  +2.83%  "
Set-TextValue "D32" "0.0548"
$ws.Range("E32").Value = "  +0.93%  "
Set-TextValue "D33" "3.98"
$ws.Range("E33").Value = "  -1.60%  "
Set-TextValue "D34" "1.27"
$ws.Range("E34").Value = "  +16.90%  "
Set-TextValue "D35" "1.77"
$ws.Range("E35").Value = "  -4.52%  "
Set-TextValue "D36" "0.710"
$ws.Range("E36").Value = "  +3.74%  "
Set-TextValue "D37" "91.90"
$ws.Range("E37").Value = "  -4.24%  "
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("D39").Value = "1.315.44"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("E41").Value = "  -0.12%  "
Set-TextValue "D42" "0.964"
$ws.Range("E42").Value = "  -3.19%  "
Set-TextValue "D43" "14.23"
$ws.Range("E43").Value = "  -2.72%  "
Set-TextValue "D44" "2.23"
$ws.Range("E44").Value = "  -8.47%  "
$ws.Range("E45").Value = "  -5.44%  "
Set-TextValue "D46" "6.19"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "1.996.94"
$ws.Range("E48").Value = "  -0.75%  "
Set-TextValue "D51" "98.95"
$ws.Range("E51").Value = "  -4.22%  "

# Rows 49/50: Cronos and PaxDollar swap places (ranking positions 47/48 stay put,
# only coin/link/price/volume move), with refreshed price + volume figures.
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D49" "1.01"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0670"
$ws.Range("E50").Value = "  +7.90%  "
